$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset ("RM 232" and "SC 92").
# Deleting row 26 first shifts "SC 92" (originally row 28) up to row 27,
# so it is removed next with a single Delete() call.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Re-roll the missing-value pattern on the remaining rows to match the new seed.
$ws.Range("D5").ClearContents()
$ws.Range("F6").Value = 16.43
$ws.Range("D8").Value = -13.9
$ws.Range("F11").Value = 17.65
$ws.Range("D12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("D14").Value = -13.1
$ws.Range("F17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("C33").Value = 10.4
